# 430104.xlsx — "Act greficos y tablas web pob"
# Rename sheets, reverse the Fecha/Valor series to descending-year order,
# relabel the metadata ("Ficha técnica") sheet with lower-case field keys,
# and append two new metadata rows (observaciones / cita / source line).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Datos" -> "Data" -------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Data"

# Data used to run ascending 2000..2020 (rows 2-22); now it runs descending
# 2020..2000 while keeping each year's value attached to that year.
$years = @(2020,2019,2018,2017,2016,2015,2014,2013,2012,2011,2010,2009,2008,2007,2006,2005,2004,2003,2002,2001,2000)
$values = @(6.2,6.8,6.8,6.5,8,7.5,7.8,8.9,9.3,10,7.7,9.6,10.6,12.1,10.6,12.7,13.2,15.1,13.7,13.9,14.1)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $years[$i]
    $ws1.Cells.Item($row, 2).Value = $values[$i]
}

# --- Sheet 2: "Ficha técnica" -> "Metadata" --------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Metadata"

$ws2.Range("A2").Value = "nomindicador"
$ws2.Range("B2").Value = "Tasa de mortalidad infantil (0 a 1 año) por 1000 nacidos vivos"

$ws2.Range("A3").Value = "derecho"
$ws2.Range("B3").Value = "Salud"

$ws2.Range("A4").Value = "conindicador"
$ws2.Range("B4").Value = "Mortalidad Infantil"

$ws2.Range("A5").Value = "tipoind"
$ws2.Range("B5").Value = "Resultados"

$ws2.Range("A6").Value = "definicion"
$ws2.Range("B6").Value = "El indicador mide la cantidad anual de defunciones en menores de un año cada 1000 nacidos vivos."

$ws2.Range("A7").Value = "calculo"
$ws2.Range("B7").Value = "Para cada año calcular: (Número de defunciones en menores de un año durante el año acaecido / Número total de nacidos vivos en el año acaecido)*1000"

$ws2.Range("A8").Value = "observaciones"
$ws2.Range("B8").Value = "Sin observaciones"

$ws2.Range("A9").Value = "cita"
$ws2.Range("B9").Value = "UMAD con base en Estadísticas Vitales - MSP"

$ws2.Range("A10").Value = "Mirador DESCA - UMAD/FCS – INDDHH"
$ws2.Range("B10").Value = " "
